# "Generate Report for Handoff" - refresh the handoff report after a new
# handoff run: the newly-generated xliff files now carry a Priority value
# ("ht") and the handoff timestamps move a bit later.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Rows 8-13 correspond to the six files that were just re-handed-off:
#   1f8e0870-..., 2c31f891-..., 717a5ec3-..., 7b7f0aec-..., 7e26ef6f-..., 84fd3000-...
$rows = 8..13

foreach ($r in $rows) {
    # Priority column (E) on the per-locale sheets goes from blank to "ht".
    $wsZhCn.Range("E$r").Value2 = "ht"
    $wsDeDe.Range("E$r").Value2 = "ht"

    # Latest Handoff Datetime (H) on zh-cn moves forward to the new run time.
    $wsZhCn.Range("H$r").Value2 = "2016-08-24 14:23:13"

    # Latest HO Xliff Generate Date (Overview!G) / Latest Handoff Datetime
    # (de-de!H) both reflect the same later timestamp for this run.
    $wsOverview.Range("G$r").Value2 = "2016-08-24 14:23:18"
    $wsDeDe.Range("H$r").Value2 = "2016-08-24 14:23:18"
}
